# flow-volumetic.xlsx update
# Adds "Default From Row" / "Default To Row" rows, a "Description" column,
# and re-themes the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Structural changes: insert two new rows (4:5) and one new column (C)
# ---------------------------------------------------------------------
$ws.Rows("4:5").Insert()
$ws.Columns("C").Insert()

# Column C should be as wide as column B (which stores width=22 in xlsx)
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# ---------------------------------------------------------------------
# 2. New row 4 / row 5 content
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "Default From Row"
$ws.Range("B4").Value = 8

$ws.Range("A5").Value = "Default To Row"
$ws.Range("B5").Value = 20
$ws.Range("D5").Value = 'Use this to select the row with the default "to" conversion value when the page loads (default is 8)'

# ---------------------------------------------------------------------
# 3. "Units" header row (now row 6) -- bold text on yellow fill
# ---------------------------------------------------------------------
$ws.Range("A6").Font.Bold = $true
$ws.Range("A6").Interior.Color = 65535

# ---------------------------------------------------------------------
# 4. "Name" / "Factor" / "Description" header row (now row 7) -- bold
# ---------------------------------------------------------------------
$ws.Range("C7").Value = "Description"
$ws.Range("A7:C7").Font.Bold = $true

# ---------------------------------------------------------------------
# 5. Description column (C8:C34) -- unit descriptions, Calibri 11 style
# ---------------------------------------------------------------------
$descriptions = @(
  "cubic meter per hour",
  "barrel per day (Petroleum)",
  "barrel per day (US)",
  "barrel per day (Imperial)",
  "barrel per minute (Petroleum)",
  "barrel per minute (US)",
  "barrel per minute (Imperial)",
  "cubic foot per minute",
  "cubic foot per second",
  "gallon per hour (US)",
  "gallon per day (US)",
  "gallon per minute (Imperial)",
  "gallon per minute (US)",
  "liter per minute",
  "liter per second",
  "liter per hour",
  "liter per day",
  "milliliter per hour",
  "cubic centimeter per hour",
  "cubic centimeter per hour",
  "cubic millimeter per hour",
  "million liters per day",
  "cubic meter per day",
  "cubic meter per minute",
  "cubic meter per second",
  "million gallons per day (Imperial)",
  "million gallons per day (US)"
)

$st = $wb.Styles.Add("Normal 2")
$st.Font.Name = "Calibri"
$st.Font.Size = 11

$startRow = 8
for ($i = 0; $i -lt $descriptions.Length; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 3).Value = $descriptions[$i]
  $ws.Cells.Item($r, 3).Style = "Normal 2"
}

# ---------------------------------------------------------------------
# 6. Misc view tweaks
# ---------------------------------------------------------------------
$ws.Range("C4").Select()
